$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$data = @(
    @('Object Name', 'Type', 'Value'),
    @('table1', 'list', 'nim$$sharma$$gaurav$$kumar'),
    @('obj2', 'label', 'kk'),
    @('obj3', 'list', 'nim$$sharma$$gaurav$$kumar'),
    @('obj4', 'label', 'tt')
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($i + 1, $j + 1).Value = $row[$j]
    }
}

$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Range("E11").Select() | Out-Null
